$p = $ppt.ActivePresentation

# --- 1. Update the cached "datetimeFigureOut" date text on the slide
#        master and every slide layout (04.11.2023 -> 05.11.2023) ---
$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        if ($shp.TextFrame.TextRange.Text -eq "04.11.2023") {
            $shp.TextFrame.TextRange.Text = "05.11.2023"
        }
    }
}

$layouts = $m.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "04.11.2023") {
                $shp.TextFrame.TextRange.Text = "05.11.2023"
            }
        }
    }
}

# --- 2. Slide 4 ("CryptoDog"): split the content placeholder text into
#        two paragraphs and fix the wording ---
$s4 = $p.Slides.Item(4)
$body4 = $s4.Shapes.Item(2).TextFrame.TextRange
$para1 = "Графический менеджер сертификатов OpenPGP на основе утилиты Gpg4win и криптографический инструмент шифрования данных"
$para2 = "Предоставляет набор функционала в области криптографии и шифрования необходимый среднестатистическому пользователю."
$body4.Text = $para1 + "`r" + $para2

# --- 3. Slide 5 ("Скрины" -> "Скриншоты") ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Скриншоты"
